$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.550.95'
$ws.Range("E2").Value = '  -0.41%  '
$ws.Range("D3").Value = '3.722.93'
$ws.Range("E3").Value = '  -1.99%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.00'
$ws.Range("E5").Value = '  -1.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '164.95'
$ws.Range("E6").Value = '  -2.55%  '
$ws.Range("D7").Value = '3.720.69'
$ws.Range("E7").Value = '  -1.89%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.517'
$ws.Range("E9").Value = '  -1.91%  '
$ws.Range("E10").Value = '  -3.94%  '
$ws.Range("E11").Value = '  -1.06%  '
$ws.Range("E12").Value = '  -2.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000261'
$ws.Range("E13").Value = '  -5.90%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.85'
$ws.Range("E14").Value = '  -2.54%  '
$ws.Range("D15").Value = '4.345.41'
$ws.Range("E15").Value = '  -2.08%  '
$ws.Range("D16").Value = '3.722.34'
$ws.Range("E16").Value = '  -1.60%  '
$ws.Range("D17").Value = '67.494.82'
$ws.Range("E17").Value = '  -0.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.22'
$ws.Range("E18").Value = '  +0.67%  '
$ws.Range("E19").Value = '  -5.13%  '
$ws.Range("E20").Value = '  -0.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.60'
$ws.Range("E21").Value = '  -1.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '463.94'
$ws.Range("E22").Value = '  -1.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.698'
$ws.Range("E23").Value = '  -3.71%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.32'
$ws.Range("E24").Value = '  -1.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000136'
$ws.Range("E25").Value = '  -10.62%  '
$ws.Range("E26").Value = '  -4.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.94'
$ws.Range("E27").Value = '  -1.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.15'
$ws.Range("E28").Value = '  -0.95%  '
$ws.Range("D30").Value = '3.866.00'
$ws.Range("E30").Value = '  -2.04%  '
$ws.Range("E31").Value = '  -6.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.30'
$ws.Range("E32").Value = '  -5.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.20'
$ws.Range("E33").Value = '  -3.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.60'
$ws.Range("E34").Value = '  -3.88%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.98'
$ws.Range("E35").Value = '  -3.64%  '
$ws.Range("D36").Value = '3.673.38'
$ws.Range("E36").Value = '  -2.42%  '
$ws.Range("E37").Value = '  -5.39%  '
$ws.Range("E38").Value = '  -10.63%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.135'
$ws.Range("E39").Value = '  -2.80%  '
$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.985'
$ws.Range("E40").Value = '  -2.62%  '
$ws.Range("E41").Value = '  -3.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.304'
$ws.Range("E44").Value = '  -3.64%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.91'
$ws.Range("E45").Value = '  -3.08%  '
$ws.Range("B46").Value = 'Cosmos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.48'
$ws.Range("E46").Value = '  -3.98%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '45.16'
$ws.Range("E47").Value = '  -2.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '392.42'
$ws.Range("E48").Value = '  -3.69%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '143.32'
$ws.Range("E49").Value = '  +1.28%  '
$ws.Range("E50").Value = '  -3.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.96'
$ws.Range("E51").Value = '  -2.56%  '
